# V19. Se agrega el código en los modelos CNN-1D, el cual que permite
# registrar las métricas en cada repetición del experimento.
# This appends 5 new rows of Loss/Accuracy metrics to the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.4852971732616425, 0.8125),
    @(0.3765910267829895, 0.8999999761581421),
    @(1.663943886756897, 0.6000000238418579),
    @(1.737413287162781, 0.7749999761581421),
    @(0.4782196879386902, 0.8062499761581421)
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
